# Update countries & provincias Spain
# Applies the 7-May-2020 18:04 data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp in A1
#  - refreshes case counts for several countries
#  - re-ranks a handful of country rows whose totals crossed over
#    each other (country name + stats move together to the new row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 1 + $i).Value = $values[$i]
    }
}

# --- Timestamp -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 18:04"

# --- Pure numeric refresh (country stays on its row) -----------------
Set-Row 4   @("Estados Unidos", 1268054, 4962, 213092, 979531, 15827, 632, 75431)
Set-Row 6   @("Italia", 215858, 1401, 96276, 89624, 1311, 274, 29958)
Set-Row 36  @("Polonia", 15047, 307, 4862, 9430, 160, 22, 755)
Set-Row 58  @("Argelia", 5182, 185, 2323, 2376, 22, 7, 483)
Set-Row 70  @("Grecia", 2678, 15, 1374, 1156, 33, 1, 148)
Set-Row 167 @("Mozambique", 81, 0, 24, 57, 0, 0, 0)

# --- Re-ranked rows (name + stats move together) ----------------------
# Colombia / Republica Dominicana swap ranking (rows 45-46)
Set-Row 45 @("Republica Dominicana", 9095, 288, 2064, 6658, 133, 11, 373)
Set-Row 46 @("Colombia", 8959, 0, 2148, 6414, 123, 0, 397)

# Kazajistan / Moldavia swap ranking (rows 59-60)
Set-Row 59 @("Moldavia", 4605, 129, 1747, 2713, 237, 2, 145)
Set-Row 60 @("Kazajistan", 4530, 108, 1470, 3030, 31, 0, 30)

# Aruba / Uganda / Nepal re-rank (rows 158-160)
Set-Row 158 @("Nepal", 101, 2, 22, 79, 0, 0, 0)
Set-Row 159 @("Aruba", 101, 0, 89, 10, 4, 0, 2)
Set-Row 160 @("Uganda", 100, 0, 55, 45, 0, 0, 0)

# Seychelles / Montserrat swap ranking (rows 205-206)
Set-Row 205 @("Montserrat", 11, 0, 7, 3, 1, 0, 1)
Set-Row 206 @("Seychelles", 11, 0, 8, 3, 0, 0, 0)
